$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6976.933
$ws.Range("J19").Value = 7186
$ws.Range("L19").Value = 7186
$ws.Range("N19").Value = -7536

$ws.Range("H33").Value = 943.06665
$ws.Range("I33").Value = 1003.4286
$ws.Range("J33").Value = 98
$ws.Range("K33").Value = 1003.4286
$ws.Range("L33").Value = 98
$ws.Range("M33").Value = -774.4286
$ws.Range("N33").Value = -556

$ws.Range("H103").Value = 691.9091
$ws.Range("I103").Value = 373.33334
$ws.Range("J103").Value = 811.375
$ws.Range("K103").Value = 1120.00002
$ws.Range("L103").Value = 2434.125
$ws.Range("M103").Value = -534.0000199999999
$ws.Range("N103").Value = -3606.125

$ws.Range("H116").Value = 4777.6665
$ws.Range("I116").Value = 4333
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4333
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -891
$ws.Range("N116").Value = -11884

$ws.Range("H127").Value = 20576.4
$ws.Range("I127").Value = 2993.25
$ws.Range("K127").Value = 8979.75
$ws.Range("M127").Value = -4019.75

$ws.Range("H132").Value = 6055.909
$ws.Range("I132").Value = 4704.2573
$ws.Range("J132").Value = 11312.333
$ws.Range("K132").Value = 14112.7719
$ws.Range("L132").Value = 33936.999
$ws.Range("M132").Value = -11582.7719
$ws.Range("N132").Value = -38996.999

$ws.Range("H135").Value = 1365.4642
$ws.Range("I135").Value = 1171.3334
$ws.Range("J135").Value = 2530.25
$ws.Range("K135").Value = 10542.0006
$ws.Range("L135").Value = 22772.25
$ws.Range("M135").Value = -8007.000599999999
$ws.Range("N135").Value = -27842.25

$ws.Range("H138").Value = 2664.242
$ws.Range("I138").Value = 1273.1578
$ws.Range("J138").Value = 3278.907
$ws.Range("K138").Value = 3819.4734
$ws.Range("L138").Value = 9836.721000000001
$ws.Range("M138").Value = 1320.5266
$ws.Range("N138").Value = -20116.721

$ws.Range("H141").Value = 13554.723
$ws.Range("J141").Value = 16131.429
$ws.Range("L141").Value = 48394.287
$ws.Range("N141").Value = -58754.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4006.3953
$ws.Range("I32").Value = 3579.1667
$ws.Range("J32").Value = 21950
$ws.Range("K32").Value = 3579.1667
$ws.Range("L32").Value = 21950
$ws.Range("M32").Value = -3292.1667
$ws.Range("N32").Value = -22524

$ws.Range("H47").Value = 39500
$ws.Range("J47").Value = 39500
$ws.Range("L47").Value = 39500
$ws.Range("N47").Value = -40950

$ws.Range("H110").Value = 5860.4546
$ws.Range("I110").Value = 4213.143
$ws.Range("K110").Value = 4213.143
$ws.Range("M110").Value = -2168.143

$ws.Range("H132").Value = 2780.6365
$ws.Range("I132").Value = 2882.5625
$ws.Range("K132").Value = 8647.6875
$ws.Range("M132").Value = -6117.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 99999
$ws.Range("J6").Value = 99999
$ws.Range("L6").Value = 99999
$ws.Range("N6").Value = -100225

$ws.Range("H86").Value = 1959.6666
$ws.Range("I86").Value = 1935.4117
$ws.Range("J86").Value = 2062.75
$ws.Range("K86").Value = 1935.4117
$ws.Range("L86").Value = 2062.75
$ws.Range("M86").Value = -812.4117000000001
$ws.Range("N86").Value = -4308.75

$ws.Range("H89").Value = 1959.6666
$ws.Range("I89").Value = 1935.4117
$ws.Range("J89").Value = 2062.75
$ws.Range("K89").Value = 9677.058500000001
$ws.Range("L89").Value = 10313.75
$ws.Range("M89").Value = -4061.058500000001
$ws.Range("N89").Value = -21545.75

$ws.Range("H94").Value = 2020.8572
$ws.Range("I94").Value = 677.2
$ws.Range("K94").Value = 677.2
$ws.Range("M94").Value = -226.2

$ws.Range("H122").Value = 179999.5
$ws.Range("J122").Value = 179999.5
$ws.Range("L122").Value = 179999.5
$ws.Range("N122").Value = -189799.5

$ws.Range("H134").Value = 3529.5405
$ws.Range("I134").Value = 2434.138
$ws.Range("K134").Value = 7302.414
$ws.Range("M134").Value = -4767.414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2685.0605
$ws.Range("I31").Value = 2565.8696
$ws.Range("K31").Value = 2565.8696
$ws.Range("M31").Value = -2270.8696

$ws.Range("H34").Value = 2685.0605
$ws.Range("I34").Value = 2565.8696
$ws.Range("K34").Value = 2565.8696
$ws.Range("M34").Value = -2363.8696

$ws.Range("H86").Value = 5132542
$ws.Range("I86").Value = 7411517.5
$ws.Range("J86").Value = 4847.5
$ws.Range("K86").Value = 7411517.5
$ws.Range("L86").Value = 4847.5
$ws.Range("M86").Value = -7410394.5
$ws.Range("N86").Value = -7093.5

$ws.Range("H89").Value = 5132542
$ws.Range("I89").Value = 7411517.5
$ws.Range("J89").Value = 4847.5
$ws.Range("K89").Value = 37057587.5
$ws.Range("L89").Value = 24237.5
$ws.Range("M89").Value = -37051971.5
$ws.Range("N89").Value = -35469.5

$ws.Range("H99").Value = 8157.1113
$ws.Range("I99").Value = 7138.4707
$ws.Range("J99").Value = 9068.526
$ws.Range("K99").Value = 7138.4707
$ws.Range("L99").Value = 9068.526
$ws.Range("M99").Value = -5640.4707
$ws.Range("N99").Value = -12064.526

$ws.Range("H126").Value = 8157.1113
$ws.Range("I126").Value = 7138.4707
$ws.Range("J126").Value = 9068.526
$ws.Range("K126").Value = 21415.4121
$ws.Range("L126").Value = 27205.578
$ws.Range("M126").Value = -18945.4121
$ws.Range("N126").Value = -32145.578

$ws.Range("H134").Value = 8340.538
$ws.Range("I134").Value = 6193.8
$ws.Range("K134").Value = 18581.4
$ws.Range("M134").Value = -16046.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.941177
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 80.454544
$ws.Range("K2").Value = 168
$ws.Range("L2").Value = 482.727264
$ws.Range("M2").Value = -55
$ws.Range("N2").Value = -708.727264

$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -338
$ws.Range("N7").ClearContents()

$ws.Range("H23").Value = 1565
$ws.Range("J23").Value = 2549.4
$ws.Range("L23").Value = 7648.200000000001
$ws.Range("N23").Value = -8118.200000000001

$ws.Range("H34").Value = 321.30768
$ws.Range("I34").Value = 87.625
$ws.Range("K34").Value = 262.875
$ws.Range("M34").Value = -178.875

$ws.Range("H39").Value = 15475.625
$ws.Range("J39").Value = 16545
$ws.Range("L39").Value = 49635
$ws.Range("N39").Value = -50223

$ws.Range("H57").Value = 13279.8
$ws.Range("J57").Value = 11299.667
$ws.Range("L57").Value = 33899.001
$ws.Range("N57").Value = -35017.001

$ws.Range("H88").Value = 4002.75
$ws.Range("J88").Value = 5335
$ws.Range("L88").Value = 16005
$ws.Range("N88").Value = -16861

$ws.Range("H91").Value = 4002.75
$ws.Range("J91").Value = 5335
$ws.Range("L91").Value = 16005
$ws.Range("N91").Value = -18969

$ws.Range("H119").Value = 6944
$ws.Range("I119").Value = 2950
$ws.Range("K119").Value = 8850
$ws.Range("M119").Value = -4012

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1637.2
$ws.Range("I70").Value = 1658.4524
$ws.Range("J70").Value = 1525.625
$ws.Range("K70").Value = 1658.4524
$ws.Range("L70").Value = 1525.625
$ws.Range("M70").Value = -1388.4524
$ws.Range("N70").Value = -2065.625

$ws.Range("H73").Value = 1637.2
$ws.Range("I73").Value = 1658.4524
$ws.Range("J73").Value = 1525.625
$ws.Range("K73").Value = 1658.4524
$ws.Range("L73").Value = 1525.625
$ws.Range("M73").Value = -722.4523999999999
$ws.Range("N73").Value = -3397.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 12222
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9774

$ws.Range("H64").Value = 42287.5
$ws.Range("J64").Value = 42287.5
$ws.Range("L64").Value = 42287.5
$ws.Range("N64").Value = -42737.5

$ws.Range("H67").Value = 42287.5
$ws.Range("J67").Value = 42287.5
$ws.Range("L67").Value = 42287.5
$ws.Range("N67").Value = -43847.5

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H93").Value = 3528.389
$ws.Range("I93").Value = 3500.7058
$ws.Range("K93").Value = 3500.7058
$ws.Range("M93").Value = -2252.7058

$ws.Range("H136").Value = 7837664.5
$ws.Range("I136").Value = 10596077
$ws.Range("J136").Value = 22164
$ws.Range("K136").Value = 31788231
$ws.Range("L136").Value = 66492
$ws.Range("M136").Value = -31785681
$ws.Range("N136").Value = -71592

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 38461730
$ws.Range("I107").Value = 180.6
$ws.Range("J107").Value = 166666900
$ws.Range("K107").Value = 541.8
$ws.Range("L107").Value = 500000700
$ws.Range("M107").Value = 1378.2
$ws.Range("N107").Value = -500004540

$ws.Range("H113").Value = 624.0833
$ws.Range("I113").Value = 348.15384
$ws.Range("J113").Value = 950.1818
$ws.Range("K113").Value = 1044.46152
$ws.Range("L113").Value = 2850.5454
$ws.Range("M113").Value = 1125.53848
$ws.Range("N113").Value = -7190.5454
